$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.556.27'
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").Value = '1.686.77'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.49%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3892'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4014'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.61%  '

$ws.Range("B9").Value = 'Polygon'
$ws.Range("C9").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.483'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("B10").Value = 'BinanceUSD'
$ws.Range("C10").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.007'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08705'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.600'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.936'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001333'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.76%  '

$ws.Range("D17").Value = '1.680.71'
$ws.Range("E17").Value = '  -0.77%  '

$ws.Range("E18").Value = '  -1.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07082'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.259'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.007'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.74%  '

$ws.Range("D24").Value = '24.591.67'
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.001'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.344'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.527'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +12.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.219'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '135.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.68%  '

$ws.Range("D32").Value = '1.868.43'
$ws.Range("E32").Value = '  -0.59%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08744'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.39%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.448'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.034'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.970'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02883'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2705'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.84%  '

$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7724'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.454'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7106'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.561'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.203'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.37%  '

$ws.Range("E48").Value = '  +0.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.344'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '90.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.48%  '
